$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("I6").Value = 2
$ws.Range("K7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 4
$ws.Range("B10").Value = 1
$ws.Range("K11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("E15").Value = 5
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = 2
$ws.Range("K16").Value = 2
$ws.Range("K17").Value = 1
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("K19").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("K20").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("B22").Value = 1
$ws.Range("C23").Value = 2
$ws.Range("F24").Value = 3
$ws.Range("K24").Value = 1
$ws.Range("E25").Value = 2
$ws.Range("F26").Value = 2
$ws.Range("E27").Value = 4
$ws.Range("I27").Value = 1
$ws.Range("I28").Value = 1
$ws.Range("K28").Value = 1
$ws.Range("C29").Value = 1
$ws.Range("K30").Value = 1
$ws.Range("C31").Value = 3
$ws.Range("C32").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("B34").Value = 1
$ws.Range("K35").Value = 1
$ws.Range("C36").Value = 1
$ws.Range("B37").Value = 1
$ws.Range("K37").Value = 2
$ws.Range("C38").Value = 1
$ws.Range("F38").Value = 1
$ws.Range("I38").Value = 1
$ws.Range("K38").Value = 1
$ws.Range("C39").Value = 1
$ws.Range("E40").Value = 1
$ws.Range("C41").Value = 1
$ws.Range("G41").Value = 2
$ws.Range("B42").Value = 1
$ws.Range("E43").Value = 3
$ws.Range("F43").Value = 1

$wb.Save()
